$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (REG VAL)
$ws.Range("B3").Value = "'0.029253814"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "'0.079765536"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = 78
$ws.Range("E3").Value = 1
$ws.Range("G3").Value = 1

# Row 5 (REG TEST)
$ws.Range("B5").Value = "'0.019499771"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "'0.10570555"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = 30
$ws.Range("E5").Value = 0
